$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A, shifting all existing
# columns (and the merged-cell header groups) one column to the right.
$ws.Columns.Item(1).Insert()

# New header label for the inserted column, in the (non-hidden) header row.
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").ClearFormats()
$ws.Range("A3").Font.Bold = $true

# Fill the data rows of the new column with the match id constant.
$ws.Range("A4:A17").Value = 16
$ws.Range("A4:A17").ClearFormats()
$ws.Range("A4:A17").Font.Bold = $true

# The hidden totals row (18) gets the same value but keeps the default
# (unbolded/unbordered) style, matching the original row's formatting.
$ws.Range("A18").Value = 16
$ws.Rows.Item(18).AutoFit()

# Restore the selection to the column that was just filled in.
[void]$ws.Range("A3:A17").Select()

Write-Output "done"
